# TC_135 - Updated loading-detail calculation values and latest-build xpath
# (Physical Layout Index switched from numeric row index to "A:n" xpath-style
# text, and LoadingDetail (column F) values recomputed for the new isolator
# loop grouping.)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- LoadingDetail (column F) recomputed for the new isolator loop grouping
$ws1.Range("F11").Value = "LI800 - 0"
$ws1.Range("F12").Value = "LI800 - 0"
$ws1.Range("F13").Value = "LI800 - 0"
$ws1.Range("F14").Value = "LI800 - 0"
$ws1.Range("F9").Value  = "801 PH - 1"
$ws1.Range("F10").Value = "801 PH - 1"
$ws1.Range("F8").Value  = "801 PH - 1"

# --- Row 4 / 5 : Isolator Units + Physical Layout Index -------------------
$ws1.Range("G4").Value = 8.5
$ws1.Range("H4").Value = "A:1"
$ws1.Range("H5").Value = "A:9"

# --- Device Quantity values (column E) -------------------------------
$ws1.Range("E8").Value  = 5
$ws1.Range("E9").Value  = 7.5
$ws1.Range("E10").Value = 8.5
$ws1.Range("E11").Value = 0
$ws1.Range("E12").Value = 5
$ws1.Range("E13").Value = 10
$ws1.Range("E14").Value = 11.5

# --- Selection / view state matches latest manual run ---------------------
$ws1.Activate()
$ws1.Range("H5").Select()
